$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (random_forest)
$ws.Range("B2").Value = 2.3872421750886934
$ws.Range("C2").Value = 0.20867501530495575
$ws.Range("D2").Value = 1.7646111111111111
$ws.Range("E2").Value = 0.42486805750385392
$ws.Range("F2").Value = 0.65181903738986779
$ws.Range("G2").Value = 0.62309714375392355
$ws.Range("H2").Value = 0.57513194249614608
$ws.Range("I2").Value = 0.76687891893576465

# Row 3 (lsboost)
$ws.Range("B3").Value = 1.4720374160240981
$ws.Range("C3").Value = 0.12867459930280581
$ws.Range("D3").Value = 1.1391348715317329
$ws.Range("E3").Value = 0.1615469719962801
$ws.Range("F3").Value = 0.40192906338840451
$ws.Range("G3").Value = 0.40223688966516002
$ws.Range("H3").Value = 0.8384530280037199
$ws.Range("I3").Value = 0.93081182381316507

# Row 4 (old_model)
$ws.Range("B4").Value = 4.2128375235700704
$ws.Range("C4").Value = 0.36825502828409712
$ws.Range("D4").Value = 3.3200000000000003
$ws.Range("E4").Value = 1.3231544574828158
$ws.Range("F4").Value = 1.1502845115373916
$ws.Range("G4").Value = 1.1723163841807913
$ws.Range("H4").Value = -0.32315445748281579
$ws.Range("I4").Value = 0.83980129657442193
